$wb = $excel.ActiveWorkbook

# New handoff GUID generated for this report run (replaces 2485d73c-c865-4c5e-b9dd-0c870353b51a)
$newGuid = "830954f8-9fd3-49ce-96f5-15d7d3de2793"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-18 20:58:40"
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.d1633fd436dd871822e2d1fa9bebe567e6e2583d.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-18 20:58:35"
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.d1633fd436dd871822e2d1fa9bebe567e6e2583d.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-18 20:58:40"
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
